$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 42547.918
$ws.Range("I86").Value = 8040.2856
$ws.Range("K86").Value = 8040.2856
$ws.Range("M86").Value = -6917.2856

$ws.Range("H89").Value = 42547.918
$ws.Range("I89").Value = 8040.2856
$ws.Range("K89").Value = 40201.428
$ws.Range("M89").Value = -34585.428

$ws.Range("H92").Value = 1203.9445
$ws.Range("I92").Value = 577.1667
$ws.Range("K92").Value = 577.1667
$ws.Range("M92").Value = 670.8333

$ws.Range("H97").Value = 2838.2917
$ws.Range("J97").Value = 3593.3333
$ws.Range("L97").Value = 10779.9999
$ws.Range("N97").Value = -11771.9999

$ws.Range("H125").Value = 4079.3333
$ws.Range("J125").Value = 4759.8
$ws.Range("L125").Value = 42838.2
$ws.Range("N125").Value = -47758.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 84995
$ws.Range("J44").Value = 84995
$ws.Range("L44").Value = 84995
$ws.Range("N44").Value = -85971

$ws.Range("H45").Value = 168833.84
$ws.Range("I45").Value = 241276.78
$ws.Range("K45").Value = 241276.78
$ws.Range("M45").Value = -240899.78

$ws.Range("H74").Value = 3111.6
$ws.Range("I74").Value = 3159.4285
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 3159.4285
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -2285.4285
$ws.Range("N74").Value = -4748

$ws.Range("H77").Value = 3111.6
$ws.Range("I77").Value = 3159.4285
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 15797.1425
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -11429.1425
$ws.Range("N77").Value = -23736

$ws.Range("H97").Value = 7155858
$ws.Range("I97").Value = 6309.4736
$ws.Range("J97").Value = 22249350
$ws.Range("K97").Value = 6309.4736
$ws.Range("L97").Value = 22249350
$ws.Range("M97").Value = -5813.4736
$ws.Range("N97").Value = -22250342

$ws.Range("H102").Value = 6086.154
$ws.Range("I102").Value = 4396.15
$ws.Range("K102").Value = 4396.15
$ws.Range("M102").Value = -2774.15

$ws.Range("H122").Value = 613498.7
$ws.Range("I122").Value = 4373.0713
$ws.Range("J122").Value = 1561027.5
$ws.Range("K122").Value = 13119.2139
$ws.Range("L122").Value = 4683082.5
$ws.Range("M122").Value = -10669.2139
$ws.Range("N122").Value = -4687982.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 72498
$ws.Range("J35").Value = 72498
$ws.Range("L35").Value = 72498
$ws.Range("N35").Value = -73118

$ws.Range("H82").Value = 64067.715
$ws.Range("J82").Value = 90618.5
$ws.Range("L82").Value = 90618.5
$ws.Range("N82").Value = -91384.5

$ws.Range("H85").Value = 64067.715
$ws.Range("J85").Value = 90618.5
$ws.Range("L85").Value = 90618.5
$ws.Range("N85").Value = -93270.5

$ws.Range("H86").Value = 4683.951
$ws.Range("I86").Value = 6235.077
$ws.Range("J86").Value = 1995.3334
$ws.Range("K86").Value = 6235.077
$ws.Range("L86").Value = 1995.3334
$ws.Range("M86").Value = -5112.077
$ws.Range("N86").Value = -4241.3334

$ws.Range("H89").Value = 4683.951
$ws.Range("I89").Value = 6235.077
$ws.Range("J89").Value = 1995.3334
$ws.Range("K89").Value = 31175.385
$ws.Range("L89").Value = 9976.666999999999
$ws.Range("M89").Value = -25559.385
$ws.Range("N89").Value = -21208.667

$ws.Range("H106").Value = 44946.332
$ws.Range("J106").Value = 44946.332
$ws.Range("L106").Value = 44946.332
$ws.Range("N106").Value = -47470.332

$ws.Range("H138").Value = 134069.22
$ws.Range("J138").Value = 134069.22
$ws.Range("L138").Value = 134069.22
$ws.Range("N138").Value = -144349.22

$ws.Range("H141").Value = 119999.664
$ws.Range("J141").Value = 119999.664
$ws.Range("L141").Value = 119999.664
$ws.Range("N141").Value = -130359.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 534.5
$ws.Range("I2").Value = 69
$ws.Range("K2").Value = 69
$ws.Range("M2").Value = 44

$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H31").Value = 3387.639
$ws.Range("I31").Value = 2398.5
$ws.Range("J31").Value = 3511.2812
$ws.Range("K31").Value = 2398.5
$ws.Range("L31").Value = 3511.2812
$ws.Range("M31").Value = -2103.5
$ws.Range("N31").Value = -4101.281199999999

$ws.Range("H34").Value = 3387.639
$ws.Range("I34").Value = 2398.5
$ws.Range("J34").Value = 3511.2812
$ws.Range("K34").Value = 2398.5
$ws.Range("L34").Value = 3511.2812
$ws.Range("M34").Value = -2196.5
$ws.Range("N34").Value = -3915.2812

$ws.Range("H58").Value = 19347.25
$ws.Range("I58").Value = 28172.8
$ws.Range("K58").Value = 28172.8
$ws.Range("M58").Value = -27969.8

$ws.Range("H86").Value = 9570
$ws.Range("I86").Value = 8879.1875
$ws.Range("J86").Value = 10798.111
$ws.Range("K86").Value = 8879.1875
$ws.Range("L86").Value = 10798.111
$ws.Range("M86").Value = -7756.1875
$ws.Range("N86").Value = -13044.111

$ws.Range("H89").Value = 9570
$ws.Range("I89").Value = 8879.1875
$ws.Range("J89").Value = 10798.111
$ws.Range("K89").Value = 44395.9375
$ws.Range("L89").Value = 53990.55500000001
$ws.Range("M89").Value = -38779.9375
$ws.Range("N89").Value = -65222.55500000001

$ws.Range("H136").Value = 19347.25
$ws.Range("I136").Value = 28172.8
$ws.Range("K136").Value = 84518.39999999999
$ws.Range("M136").Value = -81968.39999999999

$ws.Range("H140").Value = 449999
$ws.Range("J140").Value = 449999
$ws.Range("L140").Value = 449999
$ws.Range("N140").Value = -460359

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5922.8647
$ws.Range("I68").Value = 2096.6667
$ws.Range("J68").Value = 7152.7144
$ws.Range("K68").Value = 6290.000100000001
$ws.Range("L68").Value = 21458.1432
$ws.Range("M68").Value = -5479.000100000001
$ws.Range("N68").Value = -23080.1432

$ws.Range("H71").Value = 5922.8647
$ws.Range("I71").Value = 2096.6667
$ws.Range("J71").Value = 7152.7144
$ws.Range("K71").Value = 18870.0003
$ws.Range("L71").Value = 64374.4296
$ws.Range("M71").Value = -14814.0003
$ws.Range("N71").Value = -72486.4296

$ws.Range("H114").Value = 15261.429
$ws.Range("I114").Value = 420.75
$ws.Range("J114").Value = 35049
$ws.Range("K114").Value = 1262.25
$ws.Range("L114").Value = 105147
$ws.Range("M114").Value = 1991.75
$ws.Range("N114").Value = -111655

$ws.Range("H121").Value = 1251379.8
$ws.Range("I121").Value = 267
$ws.Range("K121").Value = 801
$ws.Range("M121").Value = 509

$ws.Range("H131").Value = 6107.2
$ws.Range("I131").Value = 6879.5835
$ws.Range("J131").Value = 3017.6667
$ws.Range("K131").Value = 20638.7505
$ws.Range("L131").Value = 9053.000100000001
$ws.Range("M131").Value = -15598.7505
$ws.Range("N131").Value = -19133.0001

$ws.Range("H140").Value = 2180.1667
$ws.Range("I140").Value = 2180.1667
$ws.Range("K140").Value = 6540.500100000001
$ws.Range("M140").Value = -1360.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8703.223
$ws.Range("I102").Value = 9416.125
$ws.Range("K102").Value = 9416.125
$ws.Range("M102").Value = -7794.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4330.8667
$ws.Range("I93").Value = 4330.8667
$ws.Range("K93").Value = 4330.8667
$ws.Range("M93").Value = -3082.8667

$ws.Range("H132").Value = 25501.916
$ws.Range("I132").Value = 35504.125
$ws.Range("J132").Value = 5497.5
$ws.Range("K132").Value = 106512.375
$ws.Range("L132").Value = 16492.5
$ws.Range("M132").Value = -103982.375
$ws.Range("N132").Value = -21552.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2470.9333
$ws.Range("J4").Value = 1879.25
$ws.Range("L4").Value = 1879.25
$ws.Range("N4").Value = -2105.25

$ws.Range("H62").Value = 98222.35000000001
$ws.Range("I62").Value = 143385.08
$ws.Range("J62").Value = 4422.846
$ws.Range("K62").Value = 143385.08
$ws.Range("L62").Value = 4422.846
$ws.Range("M62").Value = -142761.08
$ws.Range("N62").Value = -5670.846

$ws.Range("H65").Value = 98222.35000000001
$ws.Range("I65").Value = 143385.08
$ws.Range("J65").Value = 4422.846
$ws.Range("K65").Value = 716925.3999999999
$ws.Range("L65").Value = 22114.23
$ws.Range("M65").Value = -713805.3999999999
$ws.Range("N65").Value = -28354.23

$ws.Range("H95").Value = 172017500
$ws.Range("J95").Value = 172017500
$ws.Range("L95").Value = 172017500
$ws.Range("N95").Value = -172022992

$ws.Range("H132").Value = 20502.738
$ws.Range("I132").Value = 24535.5
$ws.Range("K132").Value = 73606.5
$ws.Range("M132").Value = -71076.5

$ws.Range("H136").Value = 4783
$ws.Range("I136").Value = 3739.6
$ws.Range("K136").Value = 11218.8
$ws.Range("M136").Value = -8668.799999999999

$ws.Range("H141").Value = 72884.75
$ws.Range("J141").Value = 72884.75
$ws.Range("L141").Value = 72884.75
$ws.Range("N141").Value = -83244.75
